# Update "想去人数" (want-to-go count) figures for two events that are
# listed on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 193   # was 191
$ws1.Range("F5").Value = 640   # was 638

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 193   # was 191
$ws4.Range("F6").Value = 640   # was 638
